$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.88
$ws.Range("J2").Value = 3.5
$ws.Range("O2").Value = 1.27
$ws.Range("P2").Value = 2.24
$ws.Range("S2").Value = 2.96
$ws.Range("T2").Value = 1.64
$ws.Range("X2").Value = 1000
$ws.Range("AA2").Value = 70
$ws.Range("AJ2").Value = 48
$ws.Range("AK2").Value = 50
$ws.Range("AL2").Value = 1000

# Row 3
$ws.Range("F3").Value = 3.7
$ws.Range("G3").Value = 4.6
$ws.Range("H3").Value = 1.8
$ws.Range("I3").Value = 1.97
$ws.Range("J3").Value = 3.7
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 1.26
$ws.Range("N3").Value = 5.7
$ws.Range("P3").Value = 2.56
$ws.Range("Q3").Value = 1.47
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 2.16
$ws.Range("T3").Value = 1.51
$ws.Range("U3").Value = 2.52
$ws.Range("V3").Value = 2.02
$ws.Range("W3").Value = 1.28
$ws.Range("Z3").Value = 16
$ws.Range("AA3").Value = 23
$ws.Range("AB3").Value = 28
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 11.5
$ws.Range("AE3").Value = 17.5
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 18.5
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 28
$ws.Range("AK3").Value = 46
$ws.Range("AL3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 7.8

# Row 4
$ws.Range("G4").Value = 12
$ws.Range("I4").Value = 1.48
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 3.05
$ws.Range("Y4").Value = 970

# Row 5
$ws.Range("H5").Value = 5.9
$ws.Range("K5").Value = 3.75
$ws.Range("N5").Value = 2.66
$ws.Range("AC5").Value = 970

# Row 6
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 1.44
$ws.Range("I6").Value = 1.5
$ws.Range("J6").Value = 4.3
$ws.Range("K6").Value = 4.8
$ws.Range("Q6").Value = 1.97
$ws.Range("T6").Value = 2.2
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 3
$ws.Range("Y6").Value = 7.2
$ws.Range("Z6").Value = 8.199999999999999
$ws.Range("AA6").Value = 13
$ws.Range("AF6").Value = 110
$ws.Range("AG6").Value = 44
$ws.Range("AH6").Value = 32
$ws.Range("AI6").Value = 50

# Row 7
$ws.Range("F7").Value = 2.96
$ws.Range("G7").Value = 3.05
$ws.Range("J7").Value = 3.05
$ws.Range("L7").Value = 1.51
$ws.Range("AA7").Value = 55
$ws.Range("AJ7").Value = 65
$ws.Range("AK7").Value = 50
$ws.Range("AN7").Value = 55

# Row 8
$ws.Range("G8").Value = 7.4
$ws.Range("P8").Value = 2.12
$ws.Range("AF8").Value = 55
$ws.Range("AI8").Value = 34

# Row 9
$ws.Range("N9").Value = 3.5
$ws.Range("P9").Value = 3.5
$ws.Range("S9").Value = 1.86
$ws.Range("W9").Value = 4.3

# Row 10
$ws.Range("I10").Value = 5.8
$ws.Range("J10").Value = 3.35
$ws.Range("K10").Value = 3.45
$ws.Range("L10").Value = 1.54
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 2.8
$ws.Range("O10").Value = 1.5
$ws.Range("P10").Value = 1.61
$ws.Range("Q10").Value = 2.5
$ws.Range("R10").Value = 1.22
$ws.Range("S10").Value = 5
$ws.Range("U10").Value = 1.74
$ws.Range("V10").Value = 1.22
$ws.Range("X10").Value = 12
$ws.Range("Y10").Value = 17
$ws.Range("Z10").Value = 42
$ws.Range("AB10").Value = 6.6
$ws.Range("AD10").Value = 23
$ws.Range("AF10").Value = 10
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 32
$ws.Range("AI10").Value = 140
$ws.Range("AK10").Value = 25
$ws.Range("AL10").Value = 55
$ws.Range("AM10").Value = 250
$ws.Range("AN10").Value = 21

# Row 11
$ws.Range("F11").Value = 1.19
$ws.Range("G11").Value = 1.22
$ws.Range("H11").Value = 15
$ws.Range("I11").Value = 17.5
$ws.Range("J11").Value = 8.6
$ws.Range("K11").Value = 10.5
$ws.Range("L11").Value = 1.17
$ws.Range("O11").Value = 1.09
$ws.Range("P11").Value = 4
$ws.Range("Q11").Value = 1.28
$ws.Range("S11").Value = 1.7
$ws.Range("T11").Value = 1.83
$ws.Range("U11").Value = 2.04
$ws.Range("V11").Value = 1.06
$ws.Range("W11").Value = 5.5
$ws.Range("X11").Value = 60
$ws.Range("Y11").Value = 80
$ws.Range("Z11").Value = 200
$ws.Range("AB11").Value = 20
$ws.Range("AC11").Value = 23
$ws.Range("AD11").Value = 60
$ws.Range("AE11").Value = 220
$ws.Range("AG11").Value = 14
$ws.Range("AI11").Value = 140
$ws.Range("AK11").Value = 14
$ws.Range("AL11").Value = 30
$ws.Range("AM11").Value = 150
$ws.Range("AN11").Value = 2.66
$ws.Range("AO11").Value = 180

# Row 12
$ws.Range("J12").Value = 4
$ws.Range("L12").Value = 1.24
$ws.Range("N12").Value = 6.8
$ws.Range("P12").Value = 2.84
$ws.Range("S12").Value = 2.04
$ws.Range("T12").Value = 1.5
$ws.Range("AA12").Value = 1000
$ws.Range("AF12").Value = 30
$ws.Range("AI12").Value = 30
$ws.Range("AJ12").Value = 50
$ws.Range("AL12").Value = 34
$ws.Range("AM12").Value = 1000
$ws.Range("AO12").Value = 11

# Row 13
$ws.Range("F13").Value = 1.5
$ws.Range("G13").Value = 1.51
$ws.Range("H13").Value = 7
$ws.Range("I13").Value = 7.2
$ws.Range("K13").Value = 5.2
$ws.Range("M13").Value = 1.04
$ws.Range("Q13").Value = 1.58
$ws.Range("R13").Value = 1.67
$ws.Range("T13").Value = 1.74
$ws.Range("W13").Value = 2.96
$ws.Range("AL13").Value = 26
$ws.Range("AM13").Value = 80
$ws.Range("AO13").Value = 70

# Row 14
$ws.Range("H14").Value = 1.95
$ws.Range("I14").Value = 2.12
$ws.Range("N14").Value = 6
$ws.Range("V14").Value = 1.89
$ws.Range("Z14").Value = 18
$ws.Range("AB14").Value = 24
$ws.Range("AE14").Value = 19.5
$ws.Range("AK14").Value = 36
$ws.Range("AN14").Value = 23
$ws.Range("AO14").Value = 8.6

# Row 15
$ws.Range("F15").Value = 2.22
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 3.35
$ws.Range("K15").Value = 3.45
$ws.Range("O15").Value = 1.39
$ws.Range("P15").Value = 1.77
$ws.Range("Q15").Value = 2.16
$ws.Range("R15").Value = 1.28
$ws.Range("V15").Value = 1.34
$ws.Range("W15").Value = 1.79
$ws.Range("X15").Value = 13.5
$ws.Range("Y15").Value = 14.5
$ws.Range("Z15").Value = 30
$ws.Range("AB15").Value = 9.4
$ws.Range("AC15").Value = 7.6
$ws.Range("AD15").Value = 17.5
$ws.Range("AF15").Value = 14.5
$ws.Range("AG15").Value = 12
$ws.Range("AH15").Value = 20
$ws.Range("AL15").Value = 55
$ws.Range("AN15").Value = 22

# Row 16
$ws.Range("F16").Value = 1.96
$ws.Range("G16").Value = 2.02
$ws.Range("H16").Value = 4.2
$ws.Range("I16").Value = 4.7
$ws.Range("J16").Value = 3.6
$ws.Range("K16").Value = 3.9
$ws.Range("N16").Value = 3.4
$ws.Range("Q16").Value = 1.92
$ws.Range("U16").Value = 1.98
$ws.Range("W16").Value = 1.98

# Row 17
$ws.Range("L17").Value = 1.16
$ws.Range("R17").Value = 2.18

# Row 18
$ws.Range("H18").Value = 4.4
$ws.Range("M18").Value = 1.08
$ws.Range("Q18").Value = 2.06
